# Populate column B ("окна пвх" / "пластиковые окна" counts) for the
# Krasnodar-region city rows. Matches the commit that adds the parsed
# per-city PVH-window / plastic-window ad counts into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{Row=8; Pvh=153; Plastik=4781},
    @{Row=9; Pvh=29; Plastik=1260},
    @{Row=10; Pvh=17; Plastik=296},
    @{Row=11; Pvh=14; Plastik=300},
    @{Row=12; Pvh=4; Plastik=257},
    @{Row=13; Pvh=11; Plastik=254},
    @{Row=14; Pvh=7; Plastik=226},
    @{Row=15; Pvh=1; Plastik=173},
    @{Row=16; Pvh=6; Plastik=265},
    @{Row=17; Pvh=3; Plastik=210},
    @{Row=18; Pvh=11; Plastik=231},
    @{Row=19; Pvh=3; Plastik=190},
    @{Row=20; Pvh=1; Plastik=121},
    @{Row=22; Pvh=3; Plastik=130},
    @{Row=23; Pvh=6; Plastik=142},
    @{Row=24; Pvh=2; Plastik=77},
    @{Row=25; Pvh=0; Plastik=64},
    @{Row=26; Pvh=2; Plastik=91},
    @{Row=27; Pvh=1; Plastik=98},
    @{Row=28; Pvh=1; Plastik=79},
    @{Row=29; Pvh=1; Plastik=63},
    @{Row=30; Pvh=0; Plastik=60},
    @{Row=31; Pvh=1; Plastik=89},
    @{Row=32; Pvh=2; Plastik=53},
    @{Row=33; Pvh=2; Plastik=48},
    @{Row=34; Pvh=0; Plastik=51},
    @{Row=35; Pvh=0; Plastik=42},
    @{Row=36; Pvh=1; Plastik=45},
    @{Row=37; Pvh=0; Plastik=31},
    @{Row=38; Pvh=0; Plastik=44},
    @{Row=39; Pvh=1; Plastik=41},
    @{Row=40; Pvh=4; Plastik=60},
    @{Row=41; Pvh=0; Plastik=20},
    @{Row=42; Pvh=0; Plastik=18},
    @{Row=43; Pvh=0; Plastik=46},
    @{Row=44; Pvh=1; Plastik=38},
    @{Row=45; Pvh=0; Plastik=24},
    @{Row=46; Pvh=0; Plastik=36},
    @{Row=47; Pvh=0; Plastik=25},
    @{Row=48; Pvh=0; Plastik=14},
    @{Row=49; Pvh=0; Plastik=26},
    @{Row=50; Pvh=0; Plastik=10},
    @{Row=51; Pvh=0; Plastik=19},
    @{Row=52; Pvh=2; Plastik=17},
    @{Row=53; Pvh=0; Plastik=4},
    @{Row=55; Pvh=1; Plastik=8},
    @{Row=56; Pvh=0; Plastik=7},
    @{Row=57; Pvh=0; Plastik=11},
    @{Row=58; Pvh=1; Plastik=8},
    @{Row=59; Pvh=0; Plastik=6},
    @{Row=60; Pvh=0; Plastik=9},
    @{Row=61; Pvh=0; Plastik=7},
    @{Row=62; Pvh=0; Plastik=6},
    @{Row=63; Pvh=0; Plastik=15},
    @{Row=64; Pvh=1; Plastik=5},
    @{Row=65; Pvh=0; Plastik=6},
    @{Row=66; Pvh=0; Plastik=6},
    @{Row=67; Pvh=0; Plastik=7},
    @{Row=68; Pvh=0; Plastik=2},
    @{Row=69; Pvh=0; Plastik=1},
    @{Row=70; Pvh=0; Plastik=3},
    @{Row=71; Pvh=0; Plastik=3},
    @{Row=72; Pvh=0; Plastik=6},
    @{Row=73; Pvh=0; Plastik=28},
    @{Row=74; Pvh=0; Plastik=2},
    @{Row=75; Pvh=2; Plastik=27},
    @{Row=76; Pvh=2; Plastik=8},
    @{Row=78; Pvh=0; Plastik=27},
    @{Row=79; Pvh=0; Plastik=3},
    @{Row=81; Pvh=0; Plastik=3},
    @{Row=82; Pvh=0; Plastik=3},
    @{Row=84; Pvh=0; Plastik=3},
    @{Row=85; Pvh=0; Plastik=1},
    @{Row=86; Pvh=0; Plastik=7},
    @{Row=87; Pvh=0; Plastik=2},
    @{Row=88; Pvh=0; Plastik=5},
    @{Row=89; Pvh=0; Plastik=5},
    @{Row=90; Pvh=0; Plastik=2},
    @{Row=91; Pvh=0; Plastik=3},
    @{Row=92; Pvh=0; Plastik=1},
    @{Row=94; Pvh=0; Plastik=2},
    @{Row=95; Pvh=0; Plastik=11},
    @{Row=96; Pvh=0; Plastik=4},
    @{Row=99; Pvh=0; Plastik=4},
    @{Row=101; Pvh=0; Plastik=2}
)

foreach ($item in $rowsData) {
    $r = $item.Row
    # remember the row height before the write so the wrap-text autofit
    # triggered by the new two-line value doesn't change row sizing
    $originalHeight = $ws.Rows.Item($r).RowHeight

    $text = "окна пвх - " + $item.Pvh + "`nпластиковые окна - " + $item.Plastik + "`n"
    $ws.Cells.Item($r, 2).Value = $text

    $ws.Rows.Item($r).RowHeight = $originalHeight
}

Write-Host ("Updated {0} rows in column B" -f $rowsData.Count)
